# modifications to TMF8801_test file
#
# 1. Delete row 5 (duplicate "APPREV_MINOR" / 0x12 entry) from the TMF8801
#    sheet - this shifts every subsequent row up by one.
# 2. Correct the "Bit Width" (col D) / "Bit Index (High)" (col E) values:
#    every register that was marked as 16-bit / bit-index-high 15 is
#    actually 8-bit / bit-index-high 7 (the 32-bit SYS_CLOCK register is
#    left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TMF8801")
$ws.Activate()

# Remove the duplicate APPREV_MINOR row.
$ws.Rows.Item(5).Delete()

# Fix up the bit-width / bit-index-high columns for the remaining rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $width = $ws.Cells.Item($r, 4).Value()
    $high  = $ws.Cells.Item($r, 5).Value()
    if ($width -eq 16) {
        $ws.Cells.Item($r, 4).Value = 8
    }
    if ($high -eq 15) {
        $ws.Cells.Item($r, 5).Value = 7
    }
}

$ws.Range("E13").Select()
